$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2606.4338
$ws.Range("I15").Value = 2606.4338
$ws.Range("K15").Value = 7819.301399999999
$ws.Range("M15").Value = -7650.301399999999

# Row 69
$ws.Range("H69").Value = 3666.3333
$ws.Range("I69").Value = 3500
$ws.Range("J69").Value = 3999
$ws.Range("K69").Value = 10500
$ws.Range("L69").Value = 11997
$ws.Range("M69").Value = -9626
$ws.Range("N69").Value = -13745

# Row 72
$ws.Range("H72").Value = 3666.3333
$ws.Range("I72").Value = 3500
$ws.Range("J72").Value = 3999
$ws.Range("K72").Value = 31500
$ws.Range("L72").Value = 35991
$ws.Range("M72").Value = -27132
$ws.Range("N72").Value = -44727

# Row 105
$ws.Range("H105").Value = 99335.5
$ws.Range("J105").Value = 99335.5
$ws.Range("L105").Value = 99335.5
$ws.Range("N105").Value = -106323.5

# Row 125
$ws.Range("H125").Value = 699.3333
$ws.Range("I125").Value = 457
$ws.Range("J125").Value = 1038.6
$ws.Range("K125").Value = 4113
$ws.Range("L125").Value = 9347.4
$ws.Range("M125").Value = -1653
$ws.Range("N125").Value = -14267.4

# Row 132
$ws.Range("H132").Value = 8156.8667
$ws.Range("I132").Value = 7704.4546
$ws.Range("J132").Value = 9401
$ws.Range("K132").Value = 23113.3638
$ws.Range("L132").Value = 28203
$ws.Range("M132").Value = -20583.3638
$ws.Range("N132").Value = -33263

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 113
$ws.Range("H113").Value = 36565.668
$ws.Range("J113").Value = 36565.668
$ws.Range("L113").Value = 36565.668
$ws.Range("N113").Value = -45243.668

# Row 132
$ws.Range("H132").Value = 4252.4585
$ws.Range("I132").Value = 4581.4375
$ws.Range("J132").Value = 3594.5
$ws.Range("K132").Value = 13744.3125
$ws.Range("L132").Value = 10783.5
$ws.Range("M132").Value = -11214.3125
$ws.Range("N132").Value = -15843.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 134
$ws.Range("H134").Value = 2686.5117
$ws.Range("I134").Value = 2274.4827
$ws.Range("J134").Value = 3540
$ws.Range("K134").Value = 6823.4481
$ws.Range("L134").Value = 10620
$ws.Range("M134").Value = -4288.4481
$ws.Range("N134").Value = -15690

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4381.553
$ws.Range("I31").Value = 927.8095
$ws.Range("J31").Value = 7171.115
$ws.Range("K31").Value = 927.8095
$ws.Range("L31").Value = 7171.115
$ws.Range("M31").Value = -632.8095
$ws.Range("N31").Value = -7761.115

# Row 34
$ws.Range("H34").Value = 4381.553
$ws.Range("I34").Value = 927.8095
$ws.Range("J34").Value = 7171.115
$ws.Range("K34").Value = 927.8095
$ws.Range("L34").Value = 7171.115
$ws.Range("M34").Value = -725.8095
$ws.Range("N34").Value = -7575.115

# Row 43
$ws.Range("H43").Value = 100000
$ws.Range("J43").Value = 100000
$ws.Range("L43").Value = 100000
$ws.Range("N43").Value = -100368

# Row 86
$ws.Range("H86").Value = 1774.9117
$ws.Range("I86").Value = 1945.8096
$ws.Range("J86").Value = 1498.8462
$ws.Range("K86").Value = 1945.8096
$ws.Range("L86").Value = 1498.8462
$ws.Range("M86").Value = -822.8096
$ws.Range("N86").Value = -3744.8462

# Row 89
$ws.Range("H89").Value = 1774.9117
$ws.Range("I89").Value = 1945.8096
$ws.Range("J89").Value = 1498.8462
$ws.Range("K89").Value = 9729.048000000001
$ws.Range("L89").Value = 7494.231
$ws.Range("M89").Value = -4113.048000000001
$ws.Range("N89").Value = -18726.231

# Row 101
$ws.Range("H101").Value = 100000
$ws.Range("J101").Value = 100000
$ws.Range("L101").Value = 100000
$ws.Range("N101").Value = -106490

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1319.1875
$ws.Range("I68").Value = 774.2857
$ws.Range("J68").Value = 1471.76
$ws.Range("K68").Value = 2322.8571
$ws.Range("L68").Value = 4415.28
$ws.Range("M68").Value = -1511.8571
$ws.Range("N68").Value = -6037.28

# Row 69
$ws.Range("H69").Value = 11111826
$ws.Range("I69").Value = 643.375
$ws.Range("J69").Value = 55556556
$ws.Range("K69").Value = 1930.125
$ws.Range("L69").Value = 166669668
$ws.Range("M69").Value = -1119.125
$ws.Range("N69").Value = -166671290

# Row 71
$ws.Range("H71").Value = 1319.1875
$ws.Range("I71").Value = 774.2857
$ws.Range("J71").Value = 1471.76
$ws.Range("K71").Value = 6968.571300000001
$ws.Range("L71").Value = 13245.84
$ws.Range("M71").Value = -2912.571300000001
$ws.Range("N71").Value = -21357.84

# Row 72
$ws.Range("H72").Value = 11111826
$ws.Range("I72").Value = 643.375
$ws.Range("J72").Value = 55556556
$ws.Range("K72").Value = 5790.375
$ws.Range("L72").Value = 500009004
$ws.Range("M72").Value = -1734.375
$ws.Range("N72").Value = -500017116

# Row 81
$ws.Range("H81").Value = 9999
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 9999
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 29997
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -32243

# Row 84
$ws.Range("H84").Value = 9999
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 9999
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 89991
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -101223

# Row 107
$ws.Range("H107").Value = 1481.6617
$ws.Range("I107").Value = 271.97144
$ws.Range("J107").Value = 2764.6667
$ws.Range("K107").Value = 815.9143199999999
$ws.Range("L107").Value = 8294.000100000001
$ws.Range("M107").Value = 1104.08568
$ws.Range("N107").Value = -12134.0001

# Row 130
$ws.Range("H130").Value = 1948.7
$ws.Range("I130").Value = 700
$ws.Range("J130").Value = 3197.4
$ws.Range("K130").Value = 2100
$ws.Range("L130").Value = 9592.200000000001
$ws.Range("M130").Value = 2920
$ws.Range("N130").Value = -19632.2

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1985.15
$ws.Range("I46").Value = 740.2
$ws.Range("J46").Value = 2400.1333
$ws.Range("K46").Value = 740.2
$ws.Range("L46").Value = 2400.1333
$ws.Range("M46").Value = -552.2
$ws.Range("N46").Value = -2776.1333

# Row 122
$ws.Range("H122").Value = 3290.2258
$ws.Range("I122").Value = 2256.682
$ws.Range("J122").Value = 5816.6665
$ws.Range("K122").Value = 6770.045999999999
$ws.Range("L122").Value = 17449.9995
$ws.Range("M122").Value = -4320.045999999999
$ws.Range("N122").Value = -22349.9995

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = 0

# Row 122
$ws.Range("H122").Value = 1861.0555
$ws.Range("I122").Value = 1524.9166
$ws.Range("K122").Value = 4574.7498
$ws.Range("M122").Value = -2124.7498

# Row 132
$ws.Range("H132").Value = 6342562.5
$ws.Range("I132").Value = 1795.3928
$ws.Range("J132").Value = 16205978
$ws.Range("K132").Value = 5386.178400000001
$ws.Range("L132").Value = 48617934
$ws.Range("M132").Value = -2856.178400000001
$ws.Range("N132").Value = -48622994

# Row 136
$ws.Range("H136").Value = 3060.1538
$ws.Range("I136").Value = 3128.6191
$ws.Range("J136").Value = 2980.2778
$ws.Range("K136").Value = 9385.8573
$ws.Range("L136").Value = 8940.8334
$ws.Range("M136").Value = -6835.8573
$ws.Range("N136").Value = -14040.8334
